$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.133.14"
$ws.Range("E2").Value = "  -3.10%  "
$ws.Range("D3").Value = "1.606.55"
$ws.Range("E3").Value = "  -2.84%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9995"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3787"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3669"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.22"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.276"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08165"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9990"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.638"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.83%  "
$ws.Range("E15").Value = "  -3.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.432"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -7.96%  "
$ws.Range("D17").Value = "1.605.16"
$ws.Range("E17").Value = "  -2.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.24"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06878"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.34"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -6.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.616"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.29%  "
$ws.Range("B22").Value = "BitDAO"
$ws.Range("C22").Value = "https://coinranking.com/coin/N2IgQ9Xme+bitdao-bit"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.5564"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.93%  "
$ws.Range("B23").Value = "Dai"
$ws.Range("C23").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("B25").Value = "WrappedBTC"
$ws.Range("C25").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D25").Value = "23.150.93"
$ws.Range("E25").Value = "  -2.98%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.372"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.46%  "
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.832"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.21%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.97%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "150.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.62%  "
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.286"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.92%  "
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "134.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("B32").Value = "WEMIXTOKEN"
$ws.Range("C32").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.385"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.23%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.931"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.24%  "
$ws.Range("B34").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C34").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D34").Value = "1.783.32"
$ws.Range("E34").Value = "  -2.72%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9639"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.28%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.07739"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.31%  "
$ws.Range("B37").Value = "FraxShare"
$ws.Range("C37").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.83%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.49%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02735"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.05%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2561"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.08901"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.371"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.53%  "
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7120"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.39%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.57%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6653"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.02%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.326"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9987"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.19%  "
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.008"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.43%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.69"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.248"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.91%  "
